# Update "想去人数" (want-to-go count) values in both the "展览" sheet
# and the "全部类型" sheet, which duplicates the same rows.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - rows 2-6, column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1046
$wsExhibit.Range("F3").Value = 266
$wsExhibit.Range("F4").Value = 2683
$wsExhibit.Range("F5").Value = 50
$wsExhibit.Range("F6").Value = 572

# Sheet "全部类型" (All Types) - rows 4-8, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1046
$wsAll.Range("F5").Value = 266
$wsAll.Range("F6").Value = 2683
$wsAll.Range("F7").Value = 50
$wsAll.Range("F8").Value = 572
